# ---------------------------------------------------------------------------
# Ping Pong Stats.xlsx - "updated data, refined add_data() function, added
# empty text file" commit replay.
#
# 1) Append 8 new games (rows 8-15) to "Ping Pong Game Data", including a new
#    player "Eric".
# 2) Re-build the F (9/11?) and I (Loser) shared formulas across the new rows.
# 3) Add a "Win %" column to "Player Records", add a row for "Eric", sort the
#    helper table by Win % descending, hide the helper rows, and add a live
#    SORT() dynamic-array table underneath.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet 1: "Ping Pong Game Data"
# =========================================================================
$ws1 = $wb.Worksheets.Item("Ping Pong Game Data")

# --- new game rows -------------------------------------------------------
$ws1.Range("A8").Value = 45319
$ws1.Range("B8").Value = "Ethan"
$ws1.Range("C8").Value = 11
$ws1.Range("D8").Value = "Eric"
$ws1.Range("E8").Value = 3
$ws1.Range("G8").Formula = "=FALSE"
$ws1.Range("H8").Formula = "=IF(C8>E8, B8, D8)"

$ws1.Range("A9").Value = 45319
$ws1.Range("B9").Value = "Andrew"
$ws1.Range("C9").Value = 11
$ws1.Range("D9").Value = "Nathan"
$ws1.Range("E9").Value = 7
$ws1.Range("G9").Formula = "=TRUE"
$ws1.Range("H9").Formula = "=IF(C9>E9, B9, D9)"

$ws1.Range("A10").Value = 45319
$ws1.Range("B10").Value = "Nathan"
$ws1.Range("C10").Value = 8
$ws1.Range("D10").Value = "Andrew"
$ws1.Range("E10").Value = 11
$ws1.Range("G10").Formula = "=FALSE"
$ws1.Range("H10").Formula = "=IF(C10>E10, B10, D10)"

$ws1.Range("A11").Value = 45319
$ws1.Range("B11").Value = "Andrew"
$ws1.Range("C11").Value = 5
$ws1.Range("D11").Value = "Nathan"
$ws1.Range("E11").Value = 11
$ws1.Range("G11").Formula = "=FALSE"
$ws1.Range("H11").Formula = "=IF(C11>E11, B11, D11)"

$ws1.Range("A12").Value = 45319
$ws1.Range("B12").Value = "Nathan"
$ws1.Range("C12").Value = 9
$ws1.Range("D12").Value = "Logan"
$ws1.Range("E12").Value = 11
$ws1.Range("G12").Formula = "=FALSE"
$ws1.Range("H12").Formula = "=IF(C12>E12, B12, D12)"

$ws1.Range("A13").Value = 45319
$ws1.Range("B13").Value = "Andrew"
$ws1.Range("C13").Value = 11
$ws1.Range("D13").Value = "Logan"
$ws1.Range("E13").Value = 5
$ws1.Range("G13").Formula = "=FALSE"
$ws1.Range("H13").Formula = "=IF(C13>E13, B13, D13)"

$ws1.Range("A14").Value = 45319
$ws1.Range("B14").Value = "Andrew"
$ws1.Range("C14").Value = 11
$ws1.Range("D14").Value = "Nathan"
$ws1.Range("E14").Value = 3
$ws1.Range("G14").Formula = "=FALSE"
$ws1.Range("H14").Formula = "=IF(C14>E14, B14, D14)"

$ws1.Range("A15").Value = 45319
$ws1.Range("B15").Value = "Logan"
$ws1.Range("C15").Value = 6
$ws1.Range("D15").Value = "Nathan"
$ws1.Range("E15").Value = 11
$ws1.Range("G15").Formula = "=FALSE"
$ws1.Range("H15").Formula = "=IF(C15>E15, B15, D15)"

# --- rebuild shared formulas ---------------------------------------------
# F2:F12 and I2:I12 grow into one shared group each (covers the pre-existing
# rows 2-7 plus the first five new rows, 8-12).
$ws1.Range("F2:F12").Formula = "=IF(AND(C2=11,E2=9),TRUE,FALSE)"
$ws1.Range("I2:I12").Formula = "=IF(C2>E2,D2,B2)"

# Rows 13-15 were filled down separately, so they form their own shared
# groups instead of extending F2:F12 / I2:I12.
$ws1.Range("F13:F15").Formula = "=IF(AND(C13=11,E13=9),TRUE,FALSE)"
$ws1.Range("I13:I15").Formula = "=IF(C13>E13,D13,B13)"

$ws1.Range("A16").Select()

# =========================================================================
# Sheet 2: "Player Records"
# =========================================================================
$ws2 = $wb.Worksheets.Item("Player Records")

# --- new "Win %" column ---------------------------------------------------
$ws2.Range("D1").Value = "Win %"

$ws2.Range("D2").Formula = "=IF((B2+C2)>0, B2/(B2+C2), 0)"
$ws2.Range("D3").Formula = "=IF((B3+C3)>0, B3/(B3+C3), 0)"
$ws2.Range("D4").Formula = "=IF((B4+C4)>0, B4/(B4+C4), 0)"
$ws2.Range("D5").Formula = "=IF((B5+C5)>0, B5/(B5+C5), 0)"
$ws2.Range("D6").Formula = "=IF((B6+C6)>0, B6/(B6+C6), 0)"
$ws2.Range("D7").Formula = "=IF((B7+C7)>0, B7/(B7+C7), 0)"
$ws2.Range("D8").Formula = "=IF((B8+C8)>0, B8/(B8+C8), 0)"
$ws2.Range("D9").Formula = "=IF((B9+C9)>0, B9/(B9+C9), 0)"

# --- new player row: Eric --------------------------------------------------
$ws2.Range("A10").Value = "Eric"
$ws2.Range("B10").Formula = "=COUNTIF('Ping Pong Game Data'!`$H`$2:`$H`$10000, A10)"
$ws2.Range("C10").Formula = "=COUNTIF('Ping Pong Game Data'!`$I`$2:`$I`$10000, A10)"
$ws2.Range("D10").Formula = "=IF((B10+C10)>0, B10/(B10+C10), 0)"

$ws2.Range("D2:D10").NumberFormat = "0.000"

# --- sort the helper table by Win % (descending) --------------------------
$sortRange = $ws2.Range("A2:D10")
$sortKey = $ws2.Range("D2:D10")
$sortRange.Sort($sortKey, 2)

$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($sortKey, 0, 2, 0, 0)
$ws2.Sort.SetRange($sortRange)
$ws2.Sort.Header = 0
$ws2.Sort.Apply()

# --- hide the helper rows, show a live sorted table below ------------------
$ws2.Rows("2:10").Hidden = $true

$ws2.Range("A11").Formula2 = "=SORT(A2:D10,4,-1)"

$ws2.Range("B1:C1").Copy()
$ws2.Range("B11:C19").PasteSpecial(-4122)
$ws2.Range("D11:D19").NumberFormat = "0.000"
$excel.CutCopyMode = 0

$ws2.Range("E12").Select()
